$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two rows after the existing data row (row 2) so the table gains two
# new example rows, exactly as a real Excel row insertion would: everything
# below shifts down by two, and the now out-of-range rows that fall past the
# sheet's last row (1048576) are removed again to keep the sheet valid.
$ws.Rows.Item(3).Resize(2).Insert()
$ws.Rows.Item(1048575).Resize(2).Delete()

# --- Row 2: first new TensorFlow example ---
$ws.Cells.Item(2, 1).Value = 36972087
$ws.Cells.Item(2, 2).Value = "2016-05-01 21:20:04"
$ws.Cells.Item(2, 3).Value = "TensorFlow"
$ws.Cells.Item(2, 4).Value = "-  train_step.run(feed_dict={x: x1[0], y_: y1[0]})`n+  train_step.run(feed_dict={x: np.expand_dims(x1[0], 0), y_: np.expand_dims(y1[0], 0)})"
$ws.Cells.Item(2, 5).Value = "Feature Input Incompatible"
$ws.Cells.Item(2, 6).Value = "Feature Data Shape"
$ws.Cells.Item(2, 7).Value = "ValueError"
$ws.Cells.Item(2, 8).Value = "Cannot feed value of shape (19,) for Tensor 'Placeholder:0', which has shape '(?, 19)'"
$ws.Rows.Item(2).RowHeight = 28.35
$ws.Cells.Item(2, 4).WrapText = $false

# --- Row 3: second new TensorFlow example ---
$ws.Cells.Item(3, 1).Value = 39009808
$ws.Cells.Item(3, 2).Value = "2016-08-18 04:03:29"
$ws.Cells.Item(3, 3).Value = "TensorFlow"
$ws.Cells.Item(3, 4).Value = "+char_num_steps = [char_num_steps]*batch_size"
$ws.Cells.Item(3, 5).Value = "Parameter Restriction Incompatible"
$ws.Cells.Item(3, 6).Value = "OP Parameter Shape"
$ws.Cells.Item(3, 7).Value = "ValueError"
$ws.Cells.Item(3, 8).Value = "sequence_length must be a vector of length batch_size, but saw shape: ()"
$ws.Rows.Item(3).RowHeight = 13.8
$ws.Cells.Item(3, 4).WrapText = $false

# --- Row 4: original Keras example, now pushed down from row 2 ---
$ws.Cells.Item(4, 1).Value = 42235611
$ws.Cells.Item(4, 2).Value = "2017-02-14 20:29:54"
$ws.Cells.Item(4, 3).Value = "Keras"
$ws.Cells.Item(4, 4).Value = "-model.add(Dense(4)) #(None, 4)`n+model.add(Dense(1))"
$ws.Cells.Item(4, 5).Value = "Label Output Incompatible"
$ws.Cells.Item(4, 6).Value = "Model Output Shape"
$ws.Cells.Item(4, 7).Value = "ValueError"
$ws.Cells.Item(4, 8).Value = "A target array with shape (100, 1) was passed for an output of shape (None, 4) while using as loss ``mean_squared_error``. This loss expects targets to have the same shape as the output."
$ws.Rows.Item(4).RowHeight = 28.35
$ws.Cells.Item(4, 4).WrapText = $true

$ws.Range("B8").Select()
